# Refresh the crypto price snapshot (Price, Volume(1h), Hora columns)
# to the latest values, keeping them stored as plain text like the
# rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:D15")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '312.58'
$ws.Range("D3").Value = '44.59'
$ws.Range("D4").Value = '5.136'
$ws.Range("D5").Value = '0.08025'
$ws.Range("D6").Value = '4.501'
$ws.Range("D7").Value = '1.653'
$ws.Range("D8").Value = '1.082'
$ws.Range("D9").Value = '0.1298'
$ws.Range("D10").Value = '0.1910'
$ws.Range("D11").Value = '0.09391'
$ws.Range("D12").Value = '0.04222'
$ws.Range("D13").Value = '0.1037'
$ws.Range("D14").Value = '0.001307'
$ws.Range("D15").Value = '0.005888'
$rng.ClearFormats()

$rng = $ws.Range("D17:D26")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '3.384'
$ws.Range("D18").Value = '2.402'
$ws.Range("D19").Value = '0.3372'
$ws.Range("D20").Value = '8.079'
$ws.Range("D21").Value = '0.1371'
$ws.Range("D22").Value = '0.3137'
$ws.Range("D23").Value = '0.04192'
$ws.Range("D24").Value = '0.001272'
$ws.Range("D25").Value = '0.004582'
$ws.Range("D26").Value = '0.0001339'
$rng.ClearFormats()

$rng = $ws.Range("D38:D51")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = '0.02655'
$ws.Range("D39").Value = '0.05411'
$ws.Range("D40").Value = '0.005623'
$ws.Range("D41").Value = '0.007708'
$ws.Range("D42").Value = '0.1412'
$ws.Range("D43").Value = '0.007330'
$ws.Range("D44").Value = '0.007864'
$ws.Range("D45").Value = '0.3122'
$ws.Range("D46").Value = '0.00006736'
$ws.Range("D47").Value = '0.00000000744'
$ws.Range("D48").Value = '0.05784'
$ws.Range("D49").Value = '0.003967'
$ws.Range("D50").Value = '0.00002083'
$ws.Range("D51").Value = '0.0001984'
$rng.ClearFormats()

$rng = $ws.Range("E2:E15")
$rng.NumberFormat = "@"
$ws.Range("E2").Value = '5.27%'
$ws.Range("E3").Value = '7.63%'
$ws.Range("E4").Value = '2.02%'
$ws.Range("E5").Value = '6.23%'
$ws.Range("E6").Value = '2.61%'
$ws.Range("E7").Value = '3.19%'
$ws.Range("E8").Value = '16.62%'
$ws.Range("E9").Value = '6.11%'
$ws.Range("E10").Value = '3.99%'
$ws.Range("E11").Value = '5.64%'
$ws.Range("E12").Value = '5.82%'
$ws.Range("E13").Value = '-1.23%'
$ws.Range("E14").Value = '0.97%'
$ws.Range("E15").Value = '-1.73%'
$rng.ClearFormats()

$rng = $ws.Range("E17:E26")
$rng.NumberFormat = "@"
$ws.Range("E17").Value = '1.16%'
$ws.Range("E18").Value = '-0.35%'
$ws.Range("E19").Value = '1.53%'
$ws.Range("E20").Value = '1.69%'
$ws.Range("E21").Value = '-3.39%'
$ws.Range("E22").Value = '4.67%'
$ws.Range("E23").Value = '3.16%'
$ws.Range("E24").Value = '0.57%'
$ws.Range("E25").Value = '15.14%'
$ws.Range("E26").Value = '8.93%'
$rng.ClearFormats()

$rng = $ws.Range("E38:E51")
$rng.NumberFormat = "@"
$ws.Range("E38").Value = '10.16%'
$ws.Range("E39").Value = '3.98%'
$ws.Range("E40").Value = '-11.94%'
$ws.Range("E41").Value = '-0.65%'
$ws.Range("E42").Value = '6.30%'
$ws.Range("E43").Value = '-1.99%'
$ws.Range("E44").Value = '0.34%'
$ws.Range("E45").Value = '-2.88%'
$ws.Range("E46").Value = '-0.58%'
$ws.Range("E47").Value = '-0.75%'
$ws.Range("E48").Value = '24.97%'
$ws.Range("E49").Value = '-5.47%'
$ws.Range("E50").Value = '-0.75%'
$ws.Range("E51").Value = '-0.75%'
$rng.ClearFormats()

$rng = $ws.Range("G2:G51")
$rng.NumberFormat = "@"
$rng.Value = '22'
$rng.ClearFormats()

